$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits -------------------------------------------
# A8: "Volume 32   Number  34" -> "Volume 32   Number  35"
$a8 = $ws.Range("A8").Value()
$idx = $a8.LastIndexOf("34") + 1
$ws.Range("A8").Characters($idx, 2).Text = "35"

# C9: "Report Covering the Week  8/18/2025  Through  8/24/2025"
#  -> "Report Covering the Week  8/25/2025  Through  8/31/2025"
$c9 = $ws.Range("C9").Value()
$idx1 = $c9.IndexOf("8/18/2025") + 1
$ws.Range("C9").Characters($idx1, 9).Text = "8/25/2025"
$c9b = $ws.Range("C9").Value()
$idx2 = $c9b.IndexOf("8/24/2025") + 1
$ws.Range("C9").Characters($idx2, 9).Text = "8/31/2025"

# --- Crime-statistics table numeric edits (rows 14-28) -------------------

# N14: {'s': '13', 't': 's', 'v': '21'} -> {'s': '15', 't': None, 'v': '0'}
$ws.Range("N14").Value = 0
$ws.Range("N14").NumberFormat = "#,##0.0;""-""#,##0.0"

# D15: {'s': '13', 't': 's', 'v': '20'} -> {'s': '14', 't': None, 'v': '1'}
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"

# E15: {'s': '13', 't': 's', 'v': '21'} -> {'s': '15', 't': None, 'v': '-100'}
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"

# G15: {'s': '13', 't': 's', 'v': '20'} -> {'s': '14', 't': None, 'v': '1'}
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = "#,##0"

# H15: {'s': '13', 't': 's', 'v': '21'} -> {'s': '15', 't': None, 'v': '-100'}
$ws.Range("H15").Value = -100
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"

# J15: {'s': '14', 't': None, 'v': '2'} -> {'s': '14', 't': None, 'v': '3'}
$ws.Range("J15").Value = 3

# K15: {'s': '15', 't': None, 'v': '150'} -> {'s': '15', 't': None, 'v': '66.666666666666'}
$ws.Range("K15").Value = 66.666666666666

# L15: {'s': '15', 't': None, 'v': '-16.666666666666'} -> {'s': '15', 't': None, 'v': '-28.571428571428'}
$ws.Range("L15").Value = -28.571428571428

# C16: {'s': '14', 't': None, 'v': '1'} -> {'s': '14', 't': None, 'v': '2'}
$ws.Range("C16").Value = 2

# D16: {'s': '14', 't': None, 'v': '4'} -> {'s': '14', 't': None, 'v': '2'}
$ws.Range("D16").Value = 2

# E16: {'s': '15', 't': None, 'v': '-75'} -> {'s': '15', 't': None, 'v': '0'}
$ws.Range("E16").Value = 0

# G16: {'s': '14', 't': None, 'v': '19'} -> {'s': '14', 't': None, 'v': '14'}
$ws.Range("G16").Value = 14

# H16: {'s': '15', 't': None, 'v': '-57.894736842105'} -> {'s': '15', 't': None, 'v': '-42.857142857142'}
$ws.Range("H16").Value = -42.857142857142

# I16: {'s': '14', 't': None, 'v': '70'} -> {'s': '14', 't': None, 'v': '72'}
$ws.Range("I16").Value = 72

# J16: {'s': '14', 't': None, 'v': '112'} -> {'s': '14', 't': None, 'v': '114'}
$ws.Range("J16").Value = 114

# K16: {'s': '15', 't': None, 'v': '-37.5'} -> {'s': '15', 't': None, 'v': '-36.842105263157'}
$ws.Range("K16").Value = -36.842105263157

# L16: {'s': '15', 't': None, 'v': '-40.677966101694'} -> {'s': '15', 't': None, 'v': '-40.495867768595'}
$ws.Range("L16").Value = -40.495867768595

# M16: {'s': '15', 't': None, 'v': '-20.454545454545'} -> {'s': '15', 't': None, 'v': '-19.101123595505'}
$ws.Range("M16").Value = -19.101123595505

# N16: {'s': '15', 't': None, 'v': '-87.387387387387'} -> {'s': '15', 't': None, 'v': '-87.256637168141'}
$ws.Range("N16").Value = -87.256637168141

# C17: {'s': '14', 't': None, 'v': '3'} -> {'s': '14', 't': None, 'v': '1'}
$ws.Range("C17").Value = 1

# E17: {'s': '15', 't': None, 'v': '-25'} -> {'s': '15', 't': None, 'v': '-75'}
$ws.Range("E17").Value = -75

# G17: {'s': '14', 't': None, 'v': '18'} -> {'s': '14', 't': None, 'v': '20'}
$ws.Range("G17").Value = 20

# H17: {'s': '15', 't': None, 'v': '-38.888888888888'} -> {'s': '15', 't': None, 'v': '-45'}
$ws.Range("H17").Value = -45

# I17: {'s': '14', 't': None, 'v': '97'} -> {'s': '14', 't': None, 'v': '98'}
$ws.Range("I17").Value = 98

# J17: {'s': '14', 't': None, 'v': '113'} -> {'s': '14', 't': None, 'v': '117'}
$ws.Range("J17").Value = 117

# K17: {'s': '15', 't': None, 'v': '-14.159292035398'} -> {'s': '15', 't': None, 'v': '-16.239316239316'}
$ws.Range("K17").Value = -16.239316239316

# L17: {'s': '15', 't': None, 'v': '-31.690140845070'} -> {'s': '15', 't': None, 'v': '-32.413793103448'}
$ws.Range("L17").Value = -32.413793103448

# M17: {'s': '15', 't': None, 'v': '64.406779661017'} -> {'s': '15', 't': None, 'v': '60.655737704918'}
$ws.Range("M17").Value = 60.655737704918

# N17: {'s': '15', 't': None, 'v': '-50.256410256410'} -> {'s': '15', 't': None, 'v': '-52.657004830917'}
$ws.Range("N17").Value = -52.657004830917

# C18: {'s': '13', 't': 's', 'v': '20'} -> {'s': '14', 't': None, 'v': '3'}
$ws.Range("C18").Value = 3
$ws.Range("C18").NumberFormat = "#,##0"

# D18: {'s': '14', 't': None, 'v': '5'} -> {'s': '14', 't': None, 'v': '3'}
$ws.Range("D18").Value = 3

# E18: {'s': '15', 't': None, 'v': '-100'} -> {'s': '15', 't': None, 'v': '0'}
$ws.Range("E18").Value = 0

# F18: {'s': '14', 't': None, 'v': '3'} -> {'s': '14', 't': None, 'v': '4'}
$ws.Range("F18").Value = 4

# G18: {'s': '14', 't': None, 'v': '16'} -> {'s': '14', 't': None, 'v': '13'}
$ws.Range("G18").Value = 13

# H18: {'s': '15', 't': None, 'v': '-81.25'} -> {'s': '15', 't': None, 'v': '-69.230769230769'}
$ws.Range("H18").Value = -69.230769230769

# I18: {'s': '14', 't': None, 'v': '103'} -> {'s': '14', 't': None, 'v': '106'}
$ws.Range("I18").Value = 106

# J18: {'s': '14', 't': None, 'v': '164'} -> {'s': '14', 't': None, 'v': '167'}
$ws.Range("J18").Value = 167

# K18: {'s': '15', 't': None, 'v': '-37.195121951219'} -> {'s': '15', 't': None, 'v': '-36.526946107784'}
$ws.Range("K18").Value = -36.526946107784

# L18: {'s': '15', 't': None, 'v': '-46.073298429319'} -> {'s': '15', 't': None, 'v': '-45.641025641025'}
$ws.Range("L18").Value = -45.641025641025

# M18: {'s': '15', 't': None, 'v': '-16.260162601626'} -> {'s': '15', 't': None, 'v': '-15.873015873015'}
$ws.Range("M18").Value = -15.873015873015

# N18: {'s': '15', 't': None, 'v': '-80.711610486891'} -> {'s': '15', 't': None, 'v': '-81.037567084078'}
$ws.Range("N18").Value = -81.037567084078

# C19: {'s': '14', 't': None, 'v': '24'} -> {'s': '14', 't': None, 'v': '25'}
$ws.Range("C19").Value = 25

# D19: {'s': '14', 't': None, 'v': '23'} -> {'s': '14', 't': None, 'v': '19'}
$ws.Range("D19").Value = 19

# E19: {'s': '15', 't': None, 'v': '4.347826086956'} -> {'s': '15', 't': None, 'v': '31.578947368421'}
$ws.Range("E19").Value = 31.578947368421

# F19: {'s': '14', 't': None, 'v': '86'} -> {'s': '14', 't': None, 'v': '94'}
$ws.Range("F19").Value = 94

# G19: {'s': '14', 't': None, 'v': '93'} -> {'s': '14', 't': None, 'v': '98'}
$ws.Range("G19").Value = 98

# H19: {'s': '15', 't': None, 'v': '-7.526881720430'} -> {'s': '15', 't': None, 'v': '-4.081632653061'}
$ws.Range("H19").Value = -4.081632653061

# I19: {'s': '14', 't': None, 'v': '631'} -> {'s': '14', 't': None, 'v': '655'}
$ws.Range("I19").Value = 655

# J19: {'s': '14', 't': None, 'v': '708'} -> {'s': '14', 't': None, 'v': '727'}
$ws.Range("J19").Value = 727

# K19: {'s': '15', 't': None, 'v': '-10.875706214689'} -> {'s': '15', 't': None, 'v': '-9.903713892709'}
$ws.Range("K19").Value = -9.903713892709

# L19: {'s': '15', 't': None, 'v': '-22.002472187886'} -> {'s': '15', 't': None, 'v': '-21.462829736211'}
$ws.Range("L19").Value = -21.462829736211

# M19: {'s': '15', 't': None, 'v': '-9.598853868194'} -> {'s': '15', 't': None, 'v': '-8.134642356241'}
$ws.Range("M19").Value = -8.134642356241

# N19: {'s': '15', 't': None, 'v': '-60.012674271229'} -> {'s': '15', 't': None, 'v': '-59.766584766584'}
$ws.Range("N19").Value = -59.766584766584

# C20: {'s': '13', 't': 's', 'v': '20'} -> {'s': '14', 't': None, 'v': '1'}
$ws.Range("C20").Value = 1
$ws.Range("C20").NumberFormat = "#,##0"

# F20: {'s': '14', 't': None, 'v': '3'} -> {'s': '14', 't': None, 'v': '2'}
$ws.Range("F20").Value = 2

# G20: {'s': '14', 't': None, 'v': '4'} -> {'s': '14', 't': None, 'v': '3'}
$ws.Range("G20").Value = 3

# H20: {'s': '15', 't': None, 'v': '-25'} -> {'s': '15', 't': None, 'v': '-33.333333333333'}
$ws.Range("H20").Value = -33.333333333333

# I20: {'s': '14', 't': None, 'v': '12'} -> {'s': '14', 't': None, 'v': '13'}
$ws.Range("I20").Value = 13

# K20: {'s': '15', 't': None, 'v': '-55.555555555555'} -> {'s': '15', 't': None, 'v': '-51.851851851851'}
$ws.Range("K20").Value = -51.851851851851

# L20: {'s': '15', 't': None, 'v': '-62.5'} -> {'s': '15', 't': None, 'v': '-59.375'}
$ws.Range("L20").Value = -59.375

# M20: {'s': '15', 't': None, 'v': '-53.846153846153'} -> {'s': '15', 't': None, 'v': '-50'}
$ws.Range("M20").Value = -50

# N20: {'s': '15', 't': None, 'v': '-97.463002114164'} -> {'s': '15', 't': None, 'v': '-97.330595482546'}
$ws.Range("N20").Value = -97.330595482546

# C21: {'s': '17', 't': None, 'v': '28'} -> {'s': '17', 't': None, 'v': '32'}
$ws.Range("C21").Value = 32

# D21: {'s': '17', 't': None, 'v': '36'} -> {'s': '17', 't': None, 'v': '29'}
$ws.Range("D21").Value = 29

# E21: {'s': '18', 't': None, 'v': '-22.222222222222'} -> {'s': '18', 't': None, 'v': '10.344827586206'}
$ws.Range("E21").Value = 10.344827586206

# F21: {'s': '17', 't': None, 'v': '111'} -> {'s': '17', 't': None, 'v': '119'}
$ws.Range("F21").Value = 119

# G21: {'s': '17', 't': None, 'v': '150'} -> {'s': '17', 't': None, 'v': '149'}
$ws.Range("G21").Value = 149

# H21: {'s': '18', 't': None, 'v': '-26'} -> {'s': '18', 't': None, 'v': '-20.134228187919'}
$ws.Range("H21").Value = -20.134228187919

# I21: {'s': '17', 't': None, 'v': '919'} -> {'s': '17', 't': None, 'v': '950'}
$ws.Range("I21").Value = 950

# J21: {'s': '17', 't': None, 'v': '1126'} -> {'s': '17', 't': None, 'v': '1155'}
$ws.Range("J21").Value = 1155

# K21: {'s': '18', 't': None, 'v': '-18.383658969804'} -> {'s': '18', 't': None, 'v': '-17.748917748917'}
$ws.Range("K21").Value = -17.748917748917

# L21: {'s': '18', 't': None, 'v': '-29.253271747498'} -> {'s': '18', 't': None, 'v': '-28.838951310861'}
$ws.Range("L21").Value = -28.838951310861

# M21: {'s': '18', 't': None, 'v': '-8.191808191808'} -> {'s': '18', 't': None, 'v': '-7.045009784735'}
$ws.Range("M21").Value = -7.045009784735

# N21: {'s': '18', 't': None, 'v': '-72.509721806760'} -> {'s': '18', 't': None, 'v': '-72.503617945007'}
$ws.Range("N21").Value = -72.503617945007

# C22: {'s': '13', 't': 's', 'v': '20'} -> {'s': '14', 't': None, 'v': '1'}
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"

# F22: {'s': '13', 't': 's', 'v': '20'} -> {'s': '14', 't': None, 'v': '1'}
$ws.Range("F22").Value = 1
$ws.Range("F22").NumberFormat = "#,##0"

# H22: {'s': '15', 't': None, 'v': '-100'} -> {'s': '15', 't': None, 'v': '-66.666666666666'}
$ws.Range("H22").Value = -66.666666666666

# I22: {'s': '14', 't': None, 'v': '30'} -> {'s': '14', 't': None, 'v': '31'}
$ws.Range("I22").Value = 31

# K22: {'s': '15', 't': None, 'v': '0'} -> {'s': '15', 't': None, 'v': '3.333333333333'}
$ws.Range("K22").Value = 3.333333333333

# L22: {'s': '15', 't': None, 'v': '3.448275862068'} -> {'s': '15', 't': None, 'v': '6.896551724137'}
$ws.Range("L22").Value = 6.896551724137

# M22: {'s': '15', 't': None, 'v': '-16.666666666666'} -> {'s': '15', 't': None, 'v': '-13.888888888888'}
$ws.Range("M22").Value = -13.888888888888

# C24: {'s': '14', 't': None, 'v': '40'} -> {'s': '14', 't': None, 'v': '28'}
$ws.Range("C24").Value = 28

# D24: {'s': '14', 't': None, 'v': '56'} -> {'s': '14', 't': None, 'v': '38'}
$ws.Range("D24").Value = 38

# E24: {'s': '15', 't': None, 'v': '-28.571428571428'} -> {'s': '15', 't': None, 'v': '-26.315789473684'}
$ws.Range("E24").Value = -26.315789473684

# F24: {'s': '14', 't': None, 'v': '157'} -> {'s': '14', 't': None, 'v': '149'}
$ws.Range("F24").Value = 149

# G24: {'s': '14', 't': None, 'v': '182'} -> {'s': '14', 't': None, 'v': '178'}
$ws.Range("G24").Value = 178

# H24: {'s': '15', 't': None, 'v': '-13.736263736263'} -> {'s': '15', 't': None, 'v': '-16.292134831460'}
$ws.Range("H24").Value = -16.292134831460

# I24: {'s': '14', 't': None, 'v': '967'} -> {'s': '14', 't': None, 'v': '993'}
$ws.Range("I24").Value = 993

# J24: {'s': '14', 't': None, 'v': '1177'} -> {'s': '14', 't': None, 'v': '1215'}
$ws.Range("J24").Value = 1215

# K24: {'s': '15', 't': None, 'v': '-17.841971112999'} -> {'s': '15', 't': None, 'v': '-18.271604938271'}
$ws.Range("K24").Value = -18.271604938271

# L24: {'s': '15', 't': None, 'v': '-29.416058394160'} -> {'s': '15', 't': None, 'v': '-29.773691654879'}
$ws.Range("L24").Value = -29.773691654879

# M24: {'s': '15', 't': None, 'v': '-2.520161290322'} -> {'s': '15', 't': None, 'v': '-2.932551319648'}
$ws.Range("M24").Value = -2.932551319648

# C25: {'s': '14', 't': None, 'v': '28'} -> {'s': '14', 't': None, 'v': '21'}
$ws.Range("C25").Value = 21

# D25: {'s': '14', 't': None, 'v': '41'} -> {'s': '14', 't': None, 'v': '23'}
$ws.Range("D25").Value = 23

# E25: {'s': '15', 't': None, 'v': '-31.707317073170'} -> {'s': '15', 't': None, 'v': '-8.695652173913'}
$ws.Range("E25").Value = -8.695652173913

# F25: {'s': '14', 't': None, 'v': '126'} -> {'s': '14', 't': None, 'v': '113'}
$ws.Range("F25").Value = 113

# G25: {'s': '14', 't': None, 'v': '139'} -> {'s': '14', 't': None, 'v': '129'}
$ws.Range("G25").Value = 129

# H25: {'s': '15', 't': None, 'v': '-9.352517985611'} -> {'s': '15', 't': None, 'v': '-12.403100775193'}
$ws.Range("H25").Value = -12.403100775193

# I25: {'s': '14', 't': None, 'v': '720'} -> {'s': '14', 't': None, 'v': '741'}
$ws.Range("I25").Value = 741

# J25: {'s': '14', 't': None, 'v': '944'} -> {'s': '14', 't': None, 'v': '967'}
$ws.Range("J25").Value = 967

# K25: {'s': '15', 't': None, 'v': '-23.728813559322'} -> {'s': '15', 't': None, 'v': '-23.371251292657'}
$ws.Range("K25").Value = -23.371251292657

# L25: {'s': '15', 't': None, 'v': '-31.297709923664'} -> {'s': '15', 't': None, 'v': '-31.515711645101'}
$ws.Range("L25").Value = -31.515711645101

# C26: {'s': '14', 't': None, 'v': '6'} -> {'s': '14', 't': None, 'v': '7'}
$ws.Range("C26").Value = 7

# D26: {'s': '14', 't': None, 'v': '6'} -> {'s': '14', 't': None, 'v': '8'}
$ws.Range("D26").Value = 8

# E26: {'s': '15', 't': None, 'v': '0'} -> {'s': '15', 't': None, 'v': '-12.5'}
$ws.Range("E26").Value = -12.5

# F26: {'s': '14', 't': None, 'v': '20'} -> {'s': '14', 't': None, 'v': '23'}
$ws.Range("F26").Value = 23

# G26: {'s': '14', 't': None, 'v': '26'} -> {'s': '14', 't': None, 'v': '23'}
$ws.Range("G26").Value = 23

# H26: {'s': '15', 't': None, 'v': '-23.076923076923'} -> {'s': '15', 't': None, 'v': '0'}
$ws.Range("H26").Value = 0

# I26: {'s': '14', 't': None, 'v': '237'} -> {'s': '14', 't': None, 'v': '244'}
$ws.Range("I26").Value = 244

# J26: {'s': '14', 't': None, 'v': '235'} -> {'s': '14', 't': None, 'v': '243'}
$ws.Range("J26").Value = 243

# K26: {'s': '15', 't': None, 'v': '0.851063829787'} -> {'s': '15', 't': None, 'v': '0.411522633744'}
$ws.Range("K26").Value = 0.411522633744

# L26: {'s': '15', 't': None, 'v': '-14.748201438848'} -> {'s': '15', 't': None, 'v': '-13.475177304964'}
$ws.Range("L26").Value = -13.475177304964

# M26: {'s': '15', 't': None, 'v': '47.204968944099'} -> {'s': '15', 't': None, 'v': '46.987951807228'}
$ws.Range("M26").Value = 46.987951807228

# D27: {'s': '13', 't': 's', 'v': '20'} -> {'s': '14', 't': None, 'v': '1'}
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"

# E27: {'s': '13', 't': 's', 'v': '21'} -> {'s': '15', 't': None, 'v': '-100'}
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"

# G27: {'s': '13', 't': 's', 'v': '20'} -> {'s': '14', 't': None, 'v': '1'}
$ws.Range("G27").Value = 1
$ws.Range("G27").NumberFormat = "#,##0"

# H27: {'s': '13', 't': 's', 'v': '21'} -> {'s': '15', 't': None, 'v': '-100'}
$ws.Range("H27").Value = -100
$ws.Range("H27").NumberFormat = "#,##0.0;""-""#,##0.0"

# J27: {'s': '14', 't': None, 'v': '6'} -> {'s': '14', 't': None, 'v': '7'}
$ws.Range("J27").Value = 7

# K27: {'s': '15', 't': None, 'v': '33.333333333333'} -> {'s': '15', 't': None, 'v': '14.285714285714'}
$ws.Range("K27").Value = 14.285714285714

# L27: {'s': '15', 't': None, 'v': '0'} -> {'s': '15', 't': None, 'v': '-11.111111111111'}
$ws.Range("L27").Value = -11.111111111111

# C28: {'s': '13', 't': 's', 'v': '20'} -> {'s': '14', 't': None, 'v': '5'}
$ws.Range("C28").Value = 5
$ws.Range("C28").NumberFormat = "#,##0"

# D28: {'s': '14', 't': None, 'v': '2'} -> {'s': '14', 't': None, 'v': '1'}
$ws.Range("D28").Value = 1

# E28: {'s': '15', 't': None, 'v': '-100'} -> {'s': '15', 't': None, 'v': '400'}
$ws.Range("E28").Value = 400

# F28: {'s': '14', 't': None, 'v': '1'} -> {'s': '14', 't': None, 'v': '5'}
$ws.Range("F28").Value = 5

# G28: {'s': '14', 't': None, 'v': '8'} -> {'s': '14', 't': None, 'v': '6'}
$ws.Range("G28").Value = 6

# H28: {'s': '15', 't': None, 'v': '-87.5'} -> {'s': '15', 't': None, 'v': '-16.666666666666'}
$ws.Range("H28").Value = -16.666666666666

# I28: {'s': '14', 't': None, 'v': '47'} -> {'s': '14', 't': None, 'v': '52'}
$ws.Range("I28").Value = 52

# J28: {'s': '14', 't': None, 'v': '48'} -> {'s': '14', 't': None, 'v': '49'}
$ws.Range("J28").Value = 49

# K28: {'s': '15', 't': None, 'v': '-2.083333333333'} -> {'s': '15', 't': None, 'v': '6.122448979591'}
$ws.Range("K28").Value = 6.122448979591

# L28: {'s': '15', 't': None, 'v': '9.302325581395'} -> {'s': '15', 't': None, 'v': '18.181818181818'}
$ws.Range("L28").Value = 18.181818181818
